$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PythonCode")
$ws.Activate()

# A7: replace the old findMaxConsecutiveOnes implementation with the new
# counter-based version (array module rewrite).
$ws.Range("A7").Value = "def findMaxConsecutiveOnes(nums) :`nmax_count = 0`ncurrent_count = 0`nfor num in nums:`nif num == 1:`ncurrent_count += 1`nmax_count = max(max_count, current_count)`n\b`n\b`nelse:`ncurrent_count = 0`n#\b\b\b`n\b`n\b`n\b`n\b`nreturn max_count"
# Assigning the new (differently-wrapped) text otherwise leaves the row
# pinned to an explicit auto-estimated height; re-autofit it so the row
# keeps using the sheet's standard height like every other row.
$ws.Rows.Item(7).AutoFit() | Out-Null

# A9: replace the old findNumbers implementation with the new one-liner.
$ws.Range("A9").Value = "def findNumbers(nums):return sum(len(str(num)) % 2 == 0 for num in nums)"

# Leave the final selection on A9, matching the saved sheet view state.
$ws.Range("A9").Select() | Out-Null
